$wb = $excel.ActiveWorkbook

$wsArch = $wb.Worksheets.Item("ARCHITECTURE")

# Update type_wall (column O) values: T2 -> T8, T6 -> T7 (new SG wall types
# based on BCA ETTV example)
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20)
foreach ($r in $rows) {
    $cell = $wsArch.Cells.Item($r, 15)  # column O
    if ($cell.Value2 -eq "T2") {
        $cell.Value = "T8"
    } elseif ($cell.Value2 -eq "T6") {
        $cell.Value = "T7"
    }
}

# Make ARCHITECTURE the active/selected sheet, with selection at O20
$wsArch.Activate()
$wsArch.Range("O20").Select()

$wb.Save()
